$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").Value = "Account"
$ws.Range("G2").Value = "m"
$ws.Range("G3").Value = "m"
$ws.Range("G4").Value = "m"
$ws.Range("G5").Value = "m"

$ws.Range("G6").Select()
